$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date (column C) from 2023-09-10 (45179) to 2023-09-11 (45180)
# for all data rows (rows 2 through 33).
for ($row = 2; $row -le 33; $row++) {
    $ws.Cells.Item($row, 3).Value = 45180
}
